$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Select the full used range on the existing (Valid_Data) sheet before
# adding the new sheet, matching the post-edit selection state there.
$ws1.Range("A1:K3").Select()

# Insert the new "Invalid_Data" sheet right after "Valid_Data" -> it
# becomes the active sheet/tab (activeTab=1, tabSelected on sheet2).
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Invalid_Data"

# ---- Header row (A1:J1) + row2/row3 (A:D, F:J) -- all reuse strings
# already present in Valid_Data, so fill order amongst these doesn't
# matter for shared-string indices.
$ws2.Range("A1").Value = "FirstName"
$ws2.Range("B1").Value = "LastName"
$ws2.Range("C1").Value = "JobTitle"
$ws2.Range("D1").Value = "CompanyName"
$ws2.Range("E1").Value = "Email"
$ws2.Range("F1").Value = "Mobile"
$ws2.Range("G1").Value = "TotalEmp"
$ws2.Range("H1").Value = "Country"
$ws2.Range("I1").Value = "HearAbout"
$ws2.Range("J1").Value = "DemoName"

$ws2.Range("A2").Value = "aaa"
$ws2.Range("B2").Value = "bbbb"
$ws2.Range("C2").Value = "ccc"
$ws2.Range("D2").Value = "ddd"

# E2 introduces the new shared string "abcd" -- it must be the FIRST
# brand-new string written so it lands at the next free shared-string
# index. Add the hyperlink (with its display text) first, then stamp
# the literal cell text; TextToDisplay sets both to "abcd@gmail.com"
# then the final .Value assignment to "abcd" is what actually mints
# the new shared string "abcd" while display="abcd@gmail.com" survives
# on the <hyperlink> element.
$ws2.Hyperlinks.Add($ws2.Range("E2"), "mailto:abcd@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "abcd@gmail.com") | Out-Null
$ws2.Range("E2").Value = "abcd"

$ws2.Range("F2").Value = 9876543210
$ws2.Range("G2").Value = 100
$ws2.Range("H2").Value = "IND"
$ws2.Range("I2").Value = "socialmedia"
$ws2.Range("J2").Value = "CompAnalyst"

$ws2.Range("A3").Value = "aab"
$ws2.Range("B3").Value = "bbc"
$ws2.Range("C3").Value = "ccd"
$ws2.Range("D3").Value = "dde"

# G3 introduces the 2nd new shared string "Test".
$ws2.Range("G3").Value = "Test"

$ws2.Hyperlinks.Add($ws2.Range("E3"), "mailto:abcd@gmail.com") | Out-Null
$ws2.Range("E3").Value = "abcd@gmail.com"

$ws2.Range("F3").Value = 9876543210
$ws2.Range("H3").Value = "IND"
$ws2.Range("I3").Value = "socialmedia"
$ws2.Range("J3").Value = "Pay Equity"

# Hyperlinks.Add leaves its own auto style behind; restore the shared
# "Hyperlink" cell style (same index used on Valid_Data!E2:E3).
$ws2.Range("E2:E3").Style = "Hyperlink"

# ---- New InvalidField / ErrorMsg columns -- order matters here too,
# to reproduce the exact shared-string indices 27..32.
$ws2.Range("K1").Value = "InvalidField"
$ws2.Range("K2").Value = "Work Email"
$ws2.Range("L1").Value = "ErrorMsg"
$ws2.Range("L2").Value = "Please input a valid email address."
$ws2.Range("K3").Value = "Total Employees"
$ws2.Range("L3").Value = "Please input a valid number."

# ---- Column widths ----
$ws2.Columns.Item(5).ColumnWidth = 22.166666666666668
$ws2.Columns.Item(6).ColumnWidth = 17.307291666666668
$ws2.Columns.Item(10).ColumnWidth = 14.451822916666666
$ws2.Columns.Item(11).ColumnWidth = 23.451822916666668
$ws2.Columns.Item(12).ColumnWidth = 40.877604166666664

# ---- Page setup ----
$ws2.PageSetup.Orientation = 1

# Leave the cursor where the original author left it.
$ws2.Range("I5").Select()
